# "Update countries & provincias Spain"
# - Re-sort a handful of country rows (Yemen, Belice/Nueva Caledonia,
#   Butan/Islas Virgenes Britanicas move to new alphabetical-ish spots)
#   which also brings a refreshed data pull for those rows.
# - Refresh daily case/death counters for several other country rows
#   that did not change position (Estados Unidos, Francia, Alemania,
#   Canada, Peru, Cuba, Etiopia).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4
$ws.Range("B4").Value = 1364956
$ws.Range("C4").Value = 17647
$ws.Range("D4").Value = 240760
$ws.Range("E4").Value = 1043499
$ws.Range("G4").Value = 660
$ws.Range("H4").Value = 80697

# Row 9
$ws.Range("B9").Value = 176867
$ws.Range("C9").Value = 209
$ws.Range("E9").Value = 94270

# Row 10
$ws.Range("B10").Value = 171780
$ws.Range("C10").Value = 456
$ws.Range("E10").Value = 19820
$ws.Range("G10").Value = 11
$ws.Range("H10").Value = 7560

# Row 15
$ws.Range("B15").Value = 68752
$ws.Range("C15").Value = 1050
$ws.Range("E15").Value = 31980

# Row 16
$ws.Range("D16").Value = 21349
$ws.Range("E16").Value = 44069

# Row 82
$ws.Range("B82").Value = 1766
$ws.Range("C82").Value = 12
$ws.Range("D82").Value = 1193
$ws.Range("E82").Value = 496
$ws.Range("F82").Value = 5
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = 77

# Row 141
$ws.Range("B141").Value = 239
$ws.Range("C141").Value = 29
$ws.Range("E141").Value = 135

# Row 174
$ws.Range("A174").Value = "Yemen"
$ws.Range("B174").Value = 51
$ws.Range("C174").Value = 17
$ws.Range("D174").Value = 1
$ws.Range("E174").Value = 42
$ws.Range("G174").Value = 1
$ws.Range("H174").Value = 8

# Row 175
$ws.Range("A175").Value = "Siria"
$ws.Range("B175").Value = 47
$ws.Range("D175").Value = 29
$ws.Range("E175").Value = 15
$ws.Range("F175").Value = 0
$ws.Range("H175").Value = 3

# Row 176
$ws.Range("A176").Value = "Macao"
$ws.Range("B176").Value = 45
$ws.Range("D176").Value = 40
$ws.Range("E176").Value = 5
$ws.Range("F176").Value = 1
$ws.Range("H176").Value = 0

# Row 177
$ws.Range("A177").Value = "Angola"
$ws.Range("B177").Value = 43
$ws.Range("D177").Value = 13
$ws.Range("H177").Value = 2

# Row 178
$ws.Range("A178").Value = "Mongolia"
$ws.Range("B178").Value = 42
$ws.Range("D178").Value = 14
$ws.Range("E178").Value = 28
$ws.Range("H178").Value = 0

# Row 179
$ws.Range("A179").Value = "Puerto Rico"
$ws.Range("D179").Value = 1
$ws.Range("E179").Value = 36
$ws.Range("F179").Value = 0
$ws.Range("H179").Value = 2

# Row 180
$ws.Range("A180").Value = "San Martin (Parte Francesa)"
$ws.Range("D180").Value = 30
$ws.Range("E180").Value = 6
$ws.Range("F180").Value = 1
$ws.Range("H180").Value = 3

# Row 181
$ws.Range("A181").Value = "Eritrea"
$ws.Range("B181").Value = 39
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 37
$ws.Range("E181").Value = 2
$ws.Range("H181").Value = 0

# Row 182
$ws.Range("A182").Value = "Zimbabue"
$ws.Range("B182").Value = 36
$ws.Range("C182").Value = 1
$ws.Range("D182").Value = 9
$ws.Range("E182").Value = 23
$ws.Range("H182").Value = 4

# Row 192
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

# Row 193
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0

# Row 212
$ws.Range("A212").Value = "Butan"
$ws.Range("D212").Value = 5
$ws.Range("H212").Value = 0

# Row 213
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 4
$ws.Range("H213").Value = 1
